# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on Sheet1 for the rows whose Strike#-derived
# value changed after recalculating against the K-based source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = 3
    4  = 4
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    12 = 1
    15 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
